$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("巨轮智能", "岩山科技", "华胜天成")
    3  = @("岩山科技", "卧龙电驱", "岩山科技")
    4  = @("阳光电源", "吉视传媒", "中际旭创")
    5  = @("浙江荣泰", "长城军工", "卧龙电驱")
    6  = @("卧龙电驱", "巨轮智能", "长城军工")
    7  = @("利欧股份", "利欧股份", "阳光电源")
    8  = @("长城军工", "阳光电源", "万通发展")
    9  = @("吉视传媒", "中际旭创", "利欧股份")
    10 = @("国轩高科", "山子高科", "新易盛")
    11 = @("中际旭创", "春兴精工", "工业富联")
    12 = @("华胜天成", "寒武纪-U", "三维通信")
    13 = @("新易盛", "华胜天成", "东方财富")
    14 = @("山子高科", "秦川机床", "北方稀土")
    15 = @("寒武纪-U", "东方财富", "方正科技")
    16 = @("通富微电", "国轩高科", "吉视传媒")
    17 = @("北方稀土", "景兴纸业", "通富微电")
    18 = @("天普股份", "浙江荣泰", "领益智造")
    19 = @("春兴精工", "新易盛", "山子高科")
    20 = @("至纯科技", "北方稀土", "天孚通信")
    21 = @("秦川机床", "通富微电", "银之杰")
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
}
